$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1807580174927114
$ws.Range("C2").Value = 0.565597667638484
$ws.Range("J2").Value = 0.02332361516034985
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.08746355685131195
$ws.Range("B3").Value = 0.005076142131979695
$ws.Range("C3").Value = 0.01015228426395939
$ws.Range("J3").Value = 0.02538071065989848
$ws.Range("P3").Value = 0.7309644670050761
$ws.Range("S3").Value = 0.2284263959390863
$ws.Range("J4").Value = 0.03636363636363636
$ws.Range("P4").Value = 0.7090909090909091
$ws.Range("S4").Value = 0.2545454545454545
$ws.Range("B6").Value = 0.07547169811320754
$ws.Range("D6").Value = 0.01509433962264151
$ws.Range("F6").Value = 0.07169811320754717
$ws.Range("J6").Value = 0.2264150943396226
$ws.Range("O6").Value = 0.03773584905660377
$ws.Range("Q6").Value = 0.1622641509433962
$ws.Range("R6").Value = 0.07924528301886792
$ws.Range("S6").Value = 0.3320754716981132
$ws.Range("B7").Value = 0.1179039301310044
$ws.Range("F7").Value = 0.03056768558951965
$ws.Range("J7").Value = 0.1441048034934498
$ws.Range("O7").Value = 0.04366812227074236
$ws.Range("Q7").Value = 0.1877729257641921
$ws.Range("R7").Value = 0.08733624454148471
$ws.Range("S7").Value = 0.388646288209607
$ws.Range("B8").Value = 0.09140767824497258
$ws.Range("D8").Value = 0.01462522851919561
$ws.Range("E8").Value = 0.005484460694698354
$ws.Range("F8").Value = 0.08226691042047532
$ws.Range("J8").Value = 0.1078610603290676
$ws.Range("O8").Value = 0.01279707495429616
$ws.Range("Q8").Value = 0.1590493601462523
$ws.Range("R8").Value = 0.129798903107861
$ws.Range("S8").Value = 0.396709323583181
$ws.Range("B9").Value = 0.06756756756756757
$ws.Range("D9").Value = 0.03378378378378379
$ws.Range("F9").Value = 0.06081081081081081
$ws.Range("J9").Value = 0.1081081081081081
$ws.Range("O9").Value = 0.03378378378378379
$ws.Range("Q9").Value = 0.1216216216216216
$ws.Range("R9").Value = 0.1756756756756757
$ws.Range("S9").Value = 0.3986486486486486
$ws.Range("B10").Value = 0.1207576953433307
$ws.Range("D10").Value = 0.03157063930544594
$ws.Range("E10").Value = 0.0007892659826361484
$ws.Range("F10").Value = 0.06866614048934491
$ws.Range("J10").Value = 0.11681136543015
$ws.Range("O10").Value = 0.0260457774269929
$ws.Range("Q10").Value = 0.1838989739542226
$ws.Range("R10").Value = 0.1112865035516969
$ws.Range("S10").Value = 0.3401736385161799
$ws.Range("G11").Value = 0.1554959785522788
$ws.Range("J11").Value = 0.08310991957104558
$ws.Range("K11").Value = 0.2144772117962467
$ws.Range("L11").Value = 0.5254691689008043
$ws.Range("S11").Value = 0.02144772117962467
$ws.Range("G12").Value = 0.751269035532995
$ws.Range("J12").Value = 0.2081218274111675
$ws.Range("K12").Value = 0.01522842639593909
$ws.Range("L12").Value = 0.01015228426395939
$ws.Range("S12").Value = 0.01522842639593909
$ws.Range("G13").Value = 0.7173913043478261
$ws.Range("J13").Value = 0.1521739130434783
$ws.Range("S13").Value = 0.1304347826086956
$ws.Range("F15").Value = 0.01626016260162602
$ws.Range("H15").Value = 0.1991869918699187
$ws.Range("I15").Value = 0.04065040650406504
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.08536585365853659
$ws.Range("M15").Value = 0.02032520325203252
$ws.Range("N15").Value = 0.004065040650406504
$ws.Range("O15").Value = 0.07317073170731707
$ws.Range("S15").Value = 0.2276422764227642
$ws.Range("F16").Value = 0.03669724770642202
$ws.Range("H16").Value = 0.1880733944954129
$ws.Range("I16").Value = 0.05045871559633028
$ws.Range("J16").Value = 0.3944954128440367
$ws.Range("K16").Value = 0.08256880733944955
$ws.Range("M16").Value = 0.02293577981651376
$ws.Range("O16").Value = 0.07339449541284404
$ws.Range("S16").Value = 0.1513761467889908
$ws.Range("F17").Value = 0.02823529411764706
$ws.Range("H17").Value = 0.2117647058823529
$ws.Range("I17").Value = 0.04941176470588235
$ws.Range("J17").Value = 0.3905882352941176
$ws.Range("K17").Value = 0.1294117647058824
$ws.Range("M17").Value = 0.02117647058823529
$ws.Range("O17").Value = 0.04941176470588235
$ws.Range("S17").Value = 0.12
$ws.Range("F18").Value = 0.02158273381294964
$ws.Range("H18").Value = 0.1798561151079137
$ws.Range("I18").Value = 0.09352517985611511
$ws.Range("J18").Value = 0.4388489208633093
$ws.Range("K18").Value = 0.08992805755395683
$ws.Range("M18").Value = 0.01079136690647482
$ws.Range("N18").Value = 0.003597122302158274
$ws.Range("O18").Value = 0.0539568345323741
$ws.Range("S18").Value = 0.1079136690647482
$ws.Range("F19").Value = 0.03255109765329296
$ws.Range("H19").Value = 0.2407267221801666
$ws.Range("I19").Value = 0.05980317940953823
$ws.Range("J19").Value = 0.3209689629068887
$ws.Range("K19").Value = 0.1180923542770628
$ws.Range("M19").Value = 0.01968205904617714
$ws.Range("O19").Value = 0.06510219530658592
$ws.Range("S19").Value = 0.1430734292202877
